$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D1 header: copy formatting from C1 (bold header style), then set text value "success"
$ws.Cells.Item(1,3).Copy()
$ws.Cells.Item(1,4).PasteSpecial(-4122)

$ws.Cells.Item(1,4).Formula = '="success"'
$ws.Cells.Item(1,4).Copy()
$ws.Cells.Item(1,4).PasteSpecial(-4163)

for ($r = 2; $r -le 144; $r++) {
    $count = $ws.Cells.Item($r, 3).Value2
    if ($count -eq 1) {
        $ws.Cells.Item($r, 4).Formula = '="0"'
    } else {
        $ws.Cells.Item($r, 4).Formula = '="1"'
    }
    $ws.Cells.Item($r, 4).Copy()
    $ws.Cells.Item($r, 4).PasteSpecial(-4163)
}
